# Append new log entries to the "Logs" worksheet, mirroring the pattern of
# existing rows (HTTP request/response log lines written by the rfid-server).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Column layout (row 1 header): A=Timestamp B=Protocol C=Method D=Route
# E=Status F=Result G=User H=UID I=RoomID J=Message K=IP

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }
$startRow = $lastRow + 1

# Each entry: Timestamp, Protocol, Method, Route, Status, Result, User, UID, RoomID, Message, IP
$rows = @(
    @("2025-11-24T21:52:47.234Z", "HTTP", "GET", "/user/6BF02F00?roomID=101", 200, "OK", "TEST", "6BF02F00", "101", "Access granted and counter incremented", "::ffff:172.28.219.204"),
    @("2025-11-24T21:52:47.272Z", "HTTP", "GET", "/uid-name/6BF02F00", 200, "OK", "TEST", "6BF02F00", "", "UID to username lookup success", "::ffff:172.28.219.204"),
    @("2025-11-24T21:52:52.655Z", "HTTP", "GET", "/user/835DF613?roomID=101", 200, "OK", "GGG", "835DF613", "101", "Access granted and counter incremented", "::ffff:172.28.219.204"),
    @("2025-11-24T21:52:52.682Z", "HTTP", "GET", "/uid-name/835DF613", 200, "OK", "GGG", "835DF613", "", "UID to username lookup success", "::ffff:172.28.219.204"),
    @("2025-11-24T21:54:22.333Z", "HTTP", "GET", "/user/835DF613?roomID=101", 200, "OK", "GGG", "835DF613", "101", "Access granted and counter incremented", "::ffff:172.28.219.204"),
    @("2025-11-24T21:54:22.372Z", "HTTP", "GET", "/uid-name/835DF613", 200, "OK", "GGG", "835DF613", "", "UID to username lookup success", "::ffff:172.28.219.204")
)

$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]

    # Column I (RoomID) looks numeric ("101") but must stay text, matching
    # every other cell in this column throughout the sheet. Enter it via a
    # text formula, then convert the formula to a plain value in place so no
    # number-format / style gets attached to the cell.
    $roomCell = $ws.Cells.Item($r, 9)
    $roomValue = [string]$row[8]
    if ($roomValue -eq "") {
        # An actual blank cell would be dropped entirely on save, but the
        # source data stores "" as a real (empty) text cell, so force text
        # entry with a quote prefix and then strip the resulting format.
        $roomCell.Value = "'"
        $roomCell.ClearFormats()
    } else {
        $roomCell.Formula = '="' + $roomValue + '"'
        $roomCell.Copy()
        $roomCell.PasteSpecial(-4163) | Out-Null
        $excel.CutCopyMode = $false
    }

    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $r = $r + 1
}

$wb.Save()
